$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Values are plain text strings (matching the sheet's existing inlineStr cells),
# so we force text via a leading apostrophe and then reset the style to Normal
# to avoid leaving a stray "quote prefix" text-format flag on the cell.

$ws.Range("D2").Value = "'305.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-4.87%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-7.56%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.026"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.36%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07670"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-6.05%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'-1.67%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.594"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-10.44%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8812"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-7.32%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.09743"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-12.50%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1720"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-7.21%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.04436"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-4.52%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08888"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-5.49%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1056"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.28%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001243"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-3.71%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005897"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.75%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.353"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.28%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.436"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-3.71%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D19").Value = "'7.038"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-5.43%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E21").Value = "'23.32%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04206"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.54%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'-4.19%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004052"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-5.74%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001222"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'9.92%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'0.01%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02314"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'0.05100"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-7.90%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007946"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.82%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1320"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-5.08%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.006500"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.61%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002023"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.98%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008688"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'2.86%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3018"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-12.99%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006554"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-6.22%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.28%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.007026"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'98.90%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003383"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-2.45%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002108"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.28%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002007"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.28%"
$ws.Range("E51").Style = "Normal"
